$wb = $excel.ActiveWorkbook

# Both the "展览" and "全部类型" sheets hold the same event listing table
# and both need their "想去人数" (F column) counts refreshed to the
# latest scraped values.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 2938
    $ws.Range("F5").Value = 6714
    $ws.Range("F6").Value = 1670
    $ws.Range("F9").Value = 56
    $ws.Range("F11").Value = 24
}
